# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.779.91"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "2.301.79"
$ws.Range("E3").Value = "  -3.16%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.33"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.60"
$ws.Range("E6").Value = "  -6.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.507"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.79"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.25"
$ws.Range("E11").Value = "  -4.14%  "
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("D15").Value = "2.658.27"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.44"
$ws.Range("D17").Value = "2.341.68"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.795"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "42.719.69"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.69"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "0.0₃0899"
$ws.Range("E22").Value = "  -4.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.35"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.95"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("E25").Value = "  -5.05%  "
$ws.Range("E26").Value = "  -4.09%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.74"
$ws.Range("E28").Value = "  -4.05%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.43"
$ws.Range("E30").Value = "  -6.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.75"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.13"
$ws.Range("E32").Value = "  -4.63%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -5.30%  "
$ws.Range("E35").Value = "  -4.80%  "
$ws.Range("E36").Value = "  -5.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.41"
$ws.Range("E37").Value = "  -6.35%  "
$ws.Range("E38").Value = "  -7.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.34"
$ws.Range("E39").Value = "  -11.23%  "
$ws.Range("E40").Value = "  -7.55%  "
$ws.Range("E41").Value = "  -5.08%  "
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("E43").Value = "  -8.68%  "
$ws.Range("D44").Value = "1.973.71"
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.30"
$ws.Range("E46").Value = "  -7.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.85"
$ws.Range("E47").Value = "  -6.84%  "
$ws.Range("E48").Value = "  -8.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.77"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "2.529.22"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.48"
$ws.Range("E51").Value = "  -7.70%  "
